$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_diagnostics")
$ws.Activate()
Write-Output $wb.ActiveSheet.Name
